$wb = $excel.ActiveWorkbook

# Insert the new "partnership" sheet right after "social care" and before "raw data"
$afterSheet = $wb.Worksheets.Item("social care")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "partnership"

# Header row
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "partnered_share"

# Data: year -> partnered_share, for 2010-2027
$years = @(2010, 2011, 2012, 2013, 2014, 2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022, 2023, 2024, 2025, 2026, 2027)
$vals  = @(0.62972779999999995, 0.62818549999999995, 0.62484090000000003, 0.61718680000000004, 0.61729670000000003, 0.61867839999999996, 0.60632839999999999, 0.60035260000000001, 0.58782400000000001, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998, 0.58603729999999998)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $vals[$i]
}

# Match the saved selection on the new sheet
$ws.Range("B1").Select()
